# risultati.xlsx: "reworked dataset and rscript, added plot pngs"
#
# On the "dataset" sheet:
#  - columns E/F ("Errate_Comprensione" / "Errate_Manutenzione") are reworked
#    into ratio columns "Rapporto_Comprensione" / "Rapporto_Manutenzione",
#    computed as Corrette_Comprensione/denominatore and
#    Corrette_Manutenzione/denominatore (denominator is the number of
#    smell/pattern items considered for that subject: 7 for rows 13-16, 8
#    for all the other rows).
#  - the W/X average columns are refreshed (same formulas/values).
#  - the saved view selection moved from S26 to T14.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("dataset")

# --- Header text: rename the two shared-string headers in place ---
$ws.Range("E1").Value = "Rapporto_Comprensione"
$ws.Range("F1").Value = "Rapporto_Manutenzione"

# --- Body rows: replace the static "Errate_*" counts with ratio formulas ---
for ($r = 2; $r -le 21; $r++) {
    $denom = 8
    if ($r -ge 13 -and $r -le 16) { $denom = 7 }
    $ws.Range("E$r").Formula = "=C$r/$denom"
    $ws.Range("F$r").Formula = "=D$r/$denom"
}

# --- Refresh the two average columns (values/formulas unchanged) ---
$ws.Range("W2:W12").Formula = "=AVERAGE(G2:N2)"
$ws.Range("X2:X12").Formula = "=AVERAGE(O2:V2)"

# --- Update the sheet's saved selection/active cell ---
$null = $ws.Range("T14").Select()
